# Update "想去人数" (want-to-go count) values per the latest scrape snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 580
$ws.Range("F3").Value = 1174
$ws.Range("F5").Value = 129
$ws.Range("F7").Value = 1246
$ws.Range("F9").Value = 94
$ws.Range("F10").Value = 933
$ws.Range("F11").Value = 908
$ws.Range("F13").Value = 84
$ws.Range("F14").Value = 90
$ws.Range("F15").Value = 659
$ws.Range("F16").Value = 867
$ws.Range("F17").Value = 1773
$ws.Range("F18").Value = 3474
$ws.Range("F19").Value = 1031
$ws.Range("F21").Value = 2429
$ws.Range("F23").Value = 33
$ws.Range("F24").Value = 3347
$ws.Range("F25").Value = 696
$ws.Range("F26").Value = 820
$ws.Range("F28").Value = 2038
$ws.Range("F29").Value = 100
$ws.Range("F30").Value = 782
$ws.Range("F33").Value = 128
$ws.Range("F35").Value = 1200
$ws.Range("F36").Value = 1860
$ws.Range("F37").Value = 452
$ws.Range("F40").Value = 234
$ws.Range("F42").Value = 211

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 16
$ws.Range("F12").Value = 105
$ws.Range("F16").Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 195

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 580
$ws.Range("F3").Value = 1174
$ws.Range("F4").Value = 129
$ws.Range("F5").Value = 1246
$ws.Range("F7").Value = 933
$ws.Range("F8").Value = 908
$ws.Range("F11").Value = 84
$ws.Range("F13").Value = 90
$ws.Range("F15").Value = 867
$ws.Range("F16").Value = 1773
$ws.Range("F17").Value = 3476
$ws.Range("F18").Value = 1031
$ws.Range("F20").Value = 2429
$ws.Range("F21").Value = 33
$ws.Range("F22").Value = 3347
$ws.Range("F23").Value = 696
$ws.Range("F24").Value = 820
$ws.Range("F27").Value = 2038
$ws.Range("F28").Value = 16
$ws.Range("F31").Value = 100
$ws.Range("F32").Value = 105
$ws.Range("F33").Value = 782
$ws.Range("F36").Value = 128
$ws.Range("F40").Value = 1200
$ws.Range("F41").Value = 1861
$ws.Range("F43").Value = 5
$ws.Range("F44").Value = 452
$ws.Range("F46").Value = 234
$ws.Range("F48").Value = 211
